$d = $word.ActiveDocument

# --- 1) Mark every inline picture as "no proofing" (<w:noProof/> on the run
#        that hosts the <w:drawing>). Two of the four pictures in this
#        document are still missing it; InlineShapes that already have it
#        are left untouched. ---
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    if ($shp.Range.NoProofing -eq 0) {
        $shp.Range.NoProofing = -1
    }
}

# --- 2) Split the caption paragraph " Elastic Net Model Learning Curve"
#        so the leading space stays in its own paragraph and the heading
#        text moves to a brand-new paragraph right after it. ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq (" Elastic Net Model Learning Curve" + [char]13)) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $spaceRange = $d.Range($target.Range.Start, $target.Range.Start + 1)
    $spaceRange.InsertParagraphAfter()
}
